# "Moved scenes to xml" -- add a new "scene" column (H) to the XML-mapped
# Sprites table, fill in its per-row scene index, and leave the selection
# on the first newly-populated cell (H9), matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the XML-mapped table by one column; Excel auto-names it "ColumnN"
# until the header cell text is set below.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# Header
$ws.Range("H1").Value2 = "scene"

# Per-row scene values (row 2..26, matches the Sprites.xml "scene" field)
$scenes = @{
    2  = 0   # Aladdin
    3  = 0   # Albert
    4  = 0   # Chinese
    5  = 0   # Tseng
    6  = 0   # Darion
    7  = 0   # Hulk
    8  = 0   # Snow White
    9  = 1   # Vadim
    10 = 2   # Barret
    11 = 2   # Bride
    12 = 2   # Captain America
    13 = 2   # Chewbacca
    14 = 2   # Dead Pool
    15 = 2   # Indiana Jones
    16 = 2   # Johhnnyyy
    17 = 2   # Rocket
    18 = 0   # Letter 1
    19 = 0   # Letter 2
    20 = 0   # Letter 3
    21 = 0   # Letter 4
    22 = 0   # Letter 5
    23 = 0   # Letter 6
    24 = 0   # Letter 7
    25 = 0   # Letter 8
    26 = 0   # Letter 9
}

foreach ($r in 2..26) {
    $ws.Cells.Item($r, 8).Value2 = $scenes[$r]
}

# New column width (best-fit-ish) for the "scene" header/data.
$ws.Columns.Item(8).ColumnWidth = 7.59

# Move the selection to where the edit was made.
$ws.Range("H9").Select() | Out-Null
